# Update test plan worksheet (Sheet1) to match the revised RDL test-case
# descriptions, re-number a few hierarchy/level values, and append two new
# rows describing the snapshot-memory vs snapshot-register test cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Level numbering tweaks -------------------------------------------------
$ws.Range("B11").Value = 1
$ws.Range("B12").Value = 1

# --- software control (sw) field descriptions: add the raw RDL token -------
$ws.Range("F13").Value = "field实例，在RDL中的sw属性为r (SW_RO)，软件只读"
$ws.Range("F14").Value = "field实例，在RDL中的sw属性为rw (SW_RW)，软件可读可写"
$ws.Range("F15").Value = "field实例，在RDL中的sw属性为rw1 (SW_RW1)，软件可读可写，但只可写一次"
$ws.Range("F16").Value = "field实例，在RDL中的sw属性为w (SW_WO)，软件只写"
$ws.Range("F17").Value = "field实例，在RDL中的sw属性为w1 (SW_W1)，软件只写，且只可写一次"

# --- alias/shared register descriptions: note differing attributes --------
$ws.Range("F28").Value = "field实例，其属于的register在RDL中声明若干（一个或多个）alias register，具有不同的hj_syncresetsignal、sw、onread、onwrite属性"
$ws.Range("F29").Value = "field实例，其属于的register在RDL中声明若干（一个或多个）shared register，具有不同的hj_syncresetsignal、sw、onread、onwrite属性"

# --- hardware control (hw) field descriptions: reorder r/rw and reword clr/set
$ws.Range("F31").Value = "field实例，在RDL中的hw属性为r (HW_RO)，硬件只读，无法通过hw_value写入数据"
$ws.Range("F32").Value = "field实例，在RDL中的hw属性为rw (HW_RW)，硬件可读可写"
$ws.Range("F33").Value = "field实例，在RDL中声明hwclr属性 (HW_CLR)，硬件访问时按位（写1）清零"
$ws.Range("F34").Value = "field实例，在RDL中声明hwset属性 (HW_SET)，硬件访问时按位（写1）置位"

# --- new "reset value" row inserted before synchronous reset, shifting the
#     remaining software-control rows down by one ---------------------------
$ws.Range("D35").Value = "reset value"
$ws.Range("F35").Value = "field实例，定义不同的reset value，观察初始值、异步复位、同步复位时是否为reset value"

$ws.Range("D36").Value = "synchronous reset"
$ws.Range("F36").Value = "field实例，在RDL中声明hj_syncresetsignal属性，带有同步复位信号入口，观察拉起后是否能够复位到reset value"

$ws.Range("D37").Value = "software modify signal"
$ws.Range("F37").Value = "field实例，在RDL中声明swmod属性，引出软件修改的指示信号，观察修改值后是否拉起"

$ws.Range("B38").Value = 7
$ws.Range("D38").Value = "software access signal"
$ws.Range("E38").Value = "DONE"
$ws.Range("F38").Value = "field实例，在RDL中声明swacc属性，引出软件访问的指示信号，观察访问后是否拉起"

$ws.Range("B39").Value = "待定"
$ws.Range("D39").Value = "write protection"
$ws.Range("E39").Value = "TBD"
$ws.Range("F39").Value = "field实例，引入写保护信号"

$ws.Range("B40").Value = 7
$ws.Range("C40").Value = "field module"
$ws.Range("D40").Value = "access precedence"
$ws.Range("F40").Value = "field实例，在RDL中声明hw或sw优先，软硬件优先级"

# --- new rows 41 & 42: split the old "snapshot register" row into separate
#     memory-access and register-access snapshot test cases ---------------
$ws.Range("A41").Value = 40
$ws.Range("B41").Value = 8
$ws.Range("C41").Value = "snapshot module"
$ws.Range("D41").Value = "snapshot memory access"
$ws.Range("E41").Value = "DONE"
$ws.Range("F41").Value = "读写请求指向外部memory，且总线数据位宽和external memory位宽不一致时，实现snapshot机制保证读写的原子性"

$ws.Range("A42").Value = 41
$ws.Range("B42").Value = 8
$ws.Range("C42").Value = "snapshot module"
$ws.Range("D42").Value = "snapshot register access"
$ws.Range("E42").Value = "TBD"
$ws.Range("F42").Value = "读写请求指向内部register，且总线数据位宽和internal register位宽不一致时，实现snapshot机制保证读写的原子性"

# --- misc view tweak (best effort; zoom level changed in the source commit)
try {
    $excel.ActiveWindow.Zoom = 100
} catch {
}
